$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row of data (row 18)
$ws.Range("A18").Value = 45960
$ws.Range("B18").Value = 647
$ws.Range("C18").Value = 10
$ws.Range("D18").Value = 637

# Update the selection to match the new active cell/range
$ws.Range("A18:D18").Select()
